$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.2199312714776632
$ws.Range("C2").Value = 0.5292096219931272
$ws.Range("J2").Value = 0.01718213058419244
$ws.Range("P2").Value = 0.140893470790378
$ws.Range("S2").Value = 0.09278350515463918

# Row 3
$ws.Range("B3").Value = 0.02597402597402598
$ws.Range("C3").Value = 0.03246753246753246
$ws.Range("J3").Value = 0.03246753246753246
$ws.Range("P3").Value = 0.7077922077922078
$ws.Range("S3").Value = 0.2012987012987013

# Row 4
$ws.Range("J4").Value = 0.07142857142857142
$ws.Range("S4").Value = 0.2619047619047619

# Row 5
$ws.Range("P5").Value = 0.3333333333333333
$ws.Range("S5").Value = 0.6666666666666666

# Row 6
$ws.Range("B6").Value = 0.04035874439461883
$ws.Range("F6").Value = 0.04035874439461883
$ws.Range("J6").Value = 0.3228699551569507
$ws.Range("Q6").Value = 0.1614349775784753
$ws.Range("R6").Value = 0.08520179372197309
$ws.Range("S6").Value = 0.3497757847533632

# Row 7
$ws.Range("B7").Value = 0.07177033492822966
$ws.Range("D7").Value = 0.02392344497607655
$ws.Range("F7").Value = 0.04784688995215311
$ws.Range("J7").Value = 0.1578947368421053
$ws.Range("O7").Value = 0.009569377990430622
$ws.Range("Q7").Value = 0.1244019138755981
$ws.Range("R7").Value = 0.138755980861244
$ws.Range("S7").Value = 0.4258373205741627

# Row 8
$ws.Range("B8").Value = 0.06593406593406594
$ws.Range("D8").Value = 0.01098901098901099
$ws.Range("F8").Value = 0.05274725274725275
$ws.Range("J8").Value = 0.1186813186813187
$ws.Range("O8").Value = 0.01758241758241758
$ws.Range("Q8").Value = 0.1626373626373626
$ws.Range("R8").Value = 0.1274725274725275
$ws.Range("S8").Value = 0.4439560439560439

# Row 9
$ws.Range("B9").Value = 0.07434944237918216
$ws.Range("D9").Value = 0.007434944237918215
$ws.Range("E9").Value = 0.007434944237918215
$ws.Range("F9").Value = 0.05947955390334572
$ws.Range("J9").Value = 0.1672862453531599
$ws.Range("O9").Value = 0.01115241635687732
$ws.Range("Q9").Value = 0.1226765799256506
$ws.Range("R9").Value = 0.1003717472118959
$ws.Range("S9").Value = 0.449814126394052

# Row 10
$ws.Range("B10").Value = 0.1094890510948905
$ws.Range("D10").Value = 0.02262773722627737
$ws.Range("E10").Value = 0.00145985401459854
$ws.Range("F10").Value = 0.06423357664233577
$ws.Range("J10").Value = 0.1262773722627737
$ws.Range("O10").Value = 0.01240875912408759
$ws.Range("Q10").Value = 0.1875912408759124
$ws.Range("R10").Value = 0.07518248175182482
$ws.Range("S10").Value = 0.4007299270072993

# Row 11
$ws.Range("G11").Value = 0.1363636363636364
$ws.Range("J11").Value = 0.07467532467532467
$ws.Range("K11").Value = 0.1753246753246753
$ws.Range("L11").Value = 0.5974025974025974
$ws.Range("S11").Value = 0.01623376623376623

# Row 12
$ws.Range("G12").Value = 0.7806122448979592
$ws.Range("J12").Value = 0.1479591836734694
$ws.Range("K12").Value = 0.00510204081632653
$ws.Range("L12").Value = 0.04591836734693878
$ws.Range("S12").Value = 0.02040816326530612

# Row 13
$ws.Range("G13").Value = 0.6176470588235294
$ws.Range("J13").Value = 0.2941176470588235
$ws.Range("S13").Value = 0.08823529411764706

# Row 15
$ws.Range("F15").Value = 0.01428571428571429
$ws.Range("H15").Value = 0.1571428571428571
$ws.Range("I15").Value = 0.09047619047619047
$ws.Range("J15").Value = 0.3857142857142857
$ws.Range("K15").Value = 0.09523809523809523
$ws.Range("M15").Value = 0.01904761904761905
$ws.Range("O15").Value = 0.0380952380952381
$ws.Range("S15").Value = 0.2

# Row 16
$ws.Range("F16").Value = 0.01714285714285714
$ws.Range("H16").Value = 0.1371428571428571
$ws.Range("I16").Value = 0.1314285714285714
$ws.Range("J16").Value = 0.4285714285714285
$ws.Range("K16").Value = 0.12
$ws.Range("M16").Value = 0.005714285714285714
$ws.Range("O16").Value = 0.06285714285714286
$ws.Range("S16").Value = 0.09714285714285714

# Row 17
$ws.Range("F17").Value = 0.02127659574468085
$ws.Range("H17").Value = 0.1843971631205674
$ws.Range("I17").Value = 0.1111111111111111
$ws.Range("J17").Value = 0.4160756501182033
$ws.Range("K17").Value = 0.09219858156028368
$ws.Range("M17").Value = 0.009456264775413711
$ws.Range("O17").Value = 0.04964539007092199
$ws.Range("S17").Value = 0.115839243498818

# Row 18
$ws.Range("F18").Value = 0.01694915254237288
$ws.Range("H18").Value = 0.1398305084745763
$ws.Range("I18").Value = 0.1398305084745763
$ws.Range("J18").Value = 0.4364406779661017
$ws.Range("K18").Value = 0.1144067796610169
$ws.Range("M18").Value = 0.01271186440677966
$ws.Range("O18").Value = 0.07627118644067797
$ws.Range("S18").Value = 0.0635593220338983

# Row 19
$ws.Range("F19").Value = 0.01783166904422254
$ws.Range("H19").Value = 0.2039942938659058
$ws.Range("I19").Value = 0.1062767475035663
$ws.Range("J19").Value = 0.3580599144079886
$ws.Range("K19").Value = 0.1034236804564907
$ws.Range("M19").Value = 0.01569186875891583
$ws.Range("O19").Value = 0.07061340941512126
$ws.Range("S19").Value = 0.1241084165477889

